# Apply the "LinuxForHealth" rebranding edit to the StructureDefinition workbook.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-job-location"

# Version bump 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date update
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher rename
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet update ---
$elements = $wb.Worksheets.Item("Elements")

# The "Extension" row's Constraint(s) cell (AI2) incorrectly duplicated the
# ele-1/ext-1 constraint text that belongs on the "Extension.extension" row
# (AI4). Clear it so it only appears once, where it is semantically correct.
$elements.Range("AI2").Value = ""

# Extension.url's Fixed Value (Q5) mirrors the same URL shared string as
# Metadata!B2; keep it in sync with the new domain.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-job-location"
